# Applies the edits described by the commit's XML diff to the Word document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...in line with module's coursework is composed by two core components.
#     Namely, illustrative Streamlit introduction and deployment on cloud..."
#    becomes
#    "...in line with module's coursework, is composed by two core
#     components. " [new paragraph] [blank paragraph] [new paragraph]
#    "Namely, Streamlit introduction and app deployment on cloud..."
# ---------------------------------------------------------------------------

# 1a) Insert a comma after "coursework".
$d.Content.Find.Execute("coursework is composed by two", $true, $false, $false, $false, $false, $true, 1, $false, "coursework, is composed by two", 2) | Out-Null

# 1b) Split "...components. Namely," into two paragraphs, with a blank
#     paragraph between them, and drop "illustrative ".
$rng = $d.Content
$rng.Find.Execute("Namely, illustrative Streamlit introduction") | Out-Null
$rng.Text = "Namely, Streamlit introduction"

$rng2 = $d.Content
$rng2.Find.Execute("Namely, Streamlit introduction") | Out-Null
$rng2.Collapse(1)
$rng2.InsertParagraphBefore()

$rng3 = $d.Content
$rng3.Find.Execute("Namely, Streamlit introduction") | Out-Null
$rng3.Collapse(1)
$rng3.InsertParagraphBefore()

# 1c) Insert "app " before "deployment on".
$d.Content.Find.Execute(" and deployment on ", $true, $false, $false, $false, $false, $true, 1, $false, " and app deployment on ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop " with" before the trailing colon.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("and specifically with:", $true, $false, $false, $false, $false, $true, 1, $false, "and specifically:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Input widgets (" -> "Play around with input widgets ("
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Input widgets (", $true, $false, $false, $false, $false, $true, 1, $false, "Play around with input widgets (", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Share your repository on " -> "Share your Streamlit app repository on "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Share your repository on ", $true, $false, $false, $false, $false, $true, 1, $false, "Share your Streamlit app repository on ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "Create an account on Heroku and GCP" -> "Create a free account on Heroku and GCP"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Create an account on Heroku and GCP", $true, $false, $false, $false, $false, $true, 1, $false, "Create a free account on Heroku and GCP", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) "Create:" -> "Create the following files:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Create:", $true, $false, $false, $false, $false, $true, 1, $false, "Create the following files:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) "Requirements.txt" -> "requirements.txt"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Requirements.txt", $true, $false, $false, $false, $false, $true, 1, $false, "requirements.txt", 2) | Out-Null

Write-Output "done"
